$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BCU")

# Delete rows 3 through 11 (keep only header row and first data row)
$ws.Range("A3:E11").EntireRow.Delete()

# Update header row
$ws.Range("A1").Value = "Qtd_Nós"
$ws.Range("B1").Value = "Ativos"
$ws.Range("C1").Value = "Distancia"
$ws.Range("D1").Value = "Tempo"
$ws.Range("E1").ClearContents()

# Update data row 2
$ws.Range("A2").Value = 81
$ws.Range("B2").Value = 35
$ws.Range("C2").Value = 11538
$ws.Range("D2").Value = 0.2519369125366211
$ws.Range("E2").ClearContents()
